$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 203.5
$ws.Range("I12").Value = 184.85715
$ws.Range("K12").Value = 184.85715
$ws.Range("M12").Value = -14.85714999999999
$ws.Range("H17").Value = 2452.1277
$ws.Range("I17").Value = 8138
$ws.Range("J17").Value = 1775.238
$ws.Range("K17").Value = 24414
$ws.Range("L17").Value = 5325.714
$ws.Range("M17").Value = -24246
$ws.Range("N17").Value = -5661.714
$ws.Range("H28").Value = 49137.816
$ws.Range("I28").Value = 64099.062
$ws.Range("J28").Value = 9241.166999999999
$ws.Range("K28").Value = 64099.062
$ws.Range("L28").Value = 9241.166999999999
$ws.Range("M28").Value = -63614.062
$ws.Range("N28").Value = -10211.167
$ws.Range("H40").Value = 3009.524
$ws.Range("I40").Value = 2749.875
$ws.Range("K40").Value = 2749.875
$ws.Range("M40").Value = -2574.875
$ws.Range("H74").Value = 3999
$ws.Range("I74").Value = 3999
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3999
$ws.Range("L74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("M74").Value = -3063
$ws.Range("H77").Value = 3999
$ws.Range("I77").Value = 3999
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 19995
$ws.Range("L77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("M77").Value = -15315
$ws.Range("H100").Value = 4302
$ws.Range("I100").Value = 3500
$ws.Range("J100").Value = 4502.5
$ws.Range("K100").Value = 3500
$ws.Range("L100").Value = 4502.5
$ws.Range("M100").Value = -2959
$ws.Range("N100").Value = -5584.5
$ws.Range("H107").Value = 533.5909
$ws.Range("I107").Value = 535.2222
$ws.Range("K107").Value = 535.2222
$ws.Range("M107").Value = 1384.7778
$ws.Range("H113").Value = 66670708
$ws.Range("I113").Value = 111114320
$ws.Range("K113").Value = 111114320
$ws.Range("M113").Value = -111111066
$ws.Range("H116").Value = 73052160
$ws.Range("I116").Value = 50208390
$ws.Range("J116").Value = 111125110
$ws.Range("K116").Value = 50208390
$ws.Range("L116").Value = 111125110
$ws.Range("M116").Value = -50204948
$ws.Range("N116").Value = -111131994
$ws.Range("H126").Value = 49149
$ws.Range("J126").Value = 49149
$ws.Range("L126").Value = 49149
$ws.Range("N126").Value = -59029
$ws.Range("H130").Value = 62245.75
$ws.Range("J130").Value = 62994.332
$ws.Range("L130").Value = 62994.332
$ws.Range("N130").Value = -73034.33199999999
$ws.Range("H132").Value = 3118.1326
$ws.Range("I132").Value = 2864.2031
$ws.Range("K132").Value = 8592.6093
$ws.Range("M132").Value = -6062.6093
$ws.Range("H137").Value = 3678.1313
$ws.Range("I137").Value = 2317
$ws.Range("J137").Value = 3848.2727
$ws.Range("K137").Value = 6951
$ws.Range("L137").Value = 11544.8181
$ws.Range("M137").Value = -4401
$ws.Range("N137").Value = -16644.8181
$ws.Range("H138").Value = 5195.711
$ws.Range("I138").Value = 3966.4
$ws.Range("J138").Value = 5546.943
$ws.Range("K138").Value = 11899.2
$ws.Range("L138").Value = 16640.829
$ws.Range("M138").Value = -6759.200000000001
$ws.Range("N138").Value = -26920.829
$ws.Range("H141").Value = 1010.2222
$ws.Range("I141").Value = 1091.0769
$ws.Range("K141").Value = 3273.2307
$ws.Range("M141").Value = 1906.7693

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1320.375
$ws.Range("I2").Value = 1223.2858
$ws.Range("K2").Value = 1223.2858
$ws.Range("M2").Value = -1110.2858
$ws.Range("H32").Value = 112505.82
$ws.Range("I32").Value = 121878.13
$ws.Range("K32").Value = 121878.13
$ws.Range("M32").Value = -121591.13
$ws.Range("H61").Value = 2426.879
$ws.Range("I61").Value = 1873.6296
$ws.Range("K61").Value = 1873.6296
$ws.Range("M61").Value = -1661.6296
$ws.Range("H74").Value = 5755.8486
$ws.Range("I74").Value = 5350.793
$ws.Range("K74").Value = 5350.793
$ws.Range("M74").Value = -4476.793
$ws.Range("H77").Value = 5755.8486
$ws.Range("I77").Value = 5350.793
$ws.Range("K77").Value = 26753.965
$ws.Range("M77").Value = -22385.965
$ws.Range("H95").Value = 13981.833
$ws.Range("J95").Value = 13981.833
$ws.Range("L95").Value = 13981.833
$ws.Range("N95").Value = -19473.833
$ws.Range("H110").Value = 71440696
$ws.Range("I110").Value = 83334820
$ws.Range("K110").Value = 83334820
$ws.Range("M110").Value = -83332775
$ws.Range("H116").Value = 1320.375
$ws.Range("I116").Value = 1223.2858
$ws.Range("K116").Value = 1223.2858
$ws.Range("M116").Value = 1070.7142
$ws.Range("H122").Value = 9525893
$ws.Range("I122").Value = 12347343
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 37042029
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -37039579
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 20835960
$ws.Range("I132").Value = 23258386
$ws.Range("K132").Value = 69775158
$ws.Range("M132").Value = -69772628
$ws.Range("H136").Value = 2426.879
$ws.Range("I136").Value = 1873.6296
$ws.Range("K136").Value = 5620.8888
$ws.Range("M136").Value = -3070.8888

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1252.7778
$ws.Range("I3").Value = 1159.375
$ws.Range("K3").Value = 1159.375
$ws.Range("M3").Value = -1045.375
$ws.Range("H13").Value = 75000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H20").Value = 10242.5
$ws.Range("I20").Value = 13648.7
$ws.Range("J20").Value = 4565.5
$ws.Range("K20").Value = 13648.7
$ws.Range("L20").Value = 4565.5
$ws.Range("M20").Value = -13401.7
$ws.Range("N20").Value = -5059.5
$ws.Range("H88").Value = 42500
$ws.Range("J88").Value = 42500
$ws.Range("L88").Value = 42500
$ws.Range("N88").Value = -43312
$ws.Range("H91").Value = 42500
$ws.Range("J91").Value = 42500
$ws.Range("L91").Value = 42500
$ws.Range("N91").Value = -45308
$ws.Range("H94").Value = 22732498
$ws.Range("I94").Value = 35718570
$ws.Range("K94").Value = 35718570
$ws.Range("M94").Value = -35718119
$ws.Range("H105").Value = 2562.25
$ws.Range("I105").Value = 2750
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 2750
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = -1003
$ws.Range("N105").Value = -5493
$ws.Range("H107").Value = 45488924
$ws.Range("I107").Value = 33439.6
$ws.Range("J107").Value = 83368500
$ws.Range("K107").Value = 33439.6
$ws.Range("L107").Value = 83368500
$ws.Range("M107").Value = -31519.6
$ws.Range("N107").Value = -83372340
$ws.Range("H109").Value = 200001
$ws.Range("J109").Value = 200001
$ws.Range("L109").Value = 200001
$ws.Range("N109").Value = -202775
$ws.Range("H134").Value = 4306.2
$ws.Range("I134").Value = 4081.182
$ws.Range("K134").Value = 12243.546
$ws.Range("M134").Value = -9708.545999999998
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 99950
$ws.Range("J9").Value = 99950
$ws.Range("L9").Value = 99950
$ws.Range("N9").Value = -100286
$ws.Range("H22").Value = 2339.3333
$ws.Range("I22").Value = 1306.1428
$ws.Range("J22").Value = 2996.818
$ws.Range("K22").Value = 1306.1428
$ws.Range("L22").Value = 2996.818
$ws.Range("M22").Value = -956.1428000000001
$ws.Range("N22").Value = -3696.818
$ws.Range("H31").Value = 3784.6667
$ws.Range("I31").Value = 880.7273
$ws.Range("J31").Value = 4308.3276
$ws.Range("K31").Value = 880.7273
$ws.Range("L31").Value = 4308.3276
$ws.Range("M31").Value = -585.7273
$ws.Range("N31").Value = -4898.3276
$ws.Range("H34").Value = 3784.6667
$ws.Range("I34").Value = 880.7273
$ws.Range("J34").Value = 4308.3276
$ws.Range("K34").Value = 880.7273
$ws.Range("L34").Value = 4308.3276
$ws.Range("M34").Value = -678.7273
$ws.Range("N34").Value = -4712.3276
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H58").Value = 457593.1
$ws.Range("I58").Value = 1616.7142
$ws.Range("J58").Value = 670382.0600000001
$ws.Range("K58").Value = 1616.7142
$ws.Range("L58").Value = 670382.0600000001
$ws.Range("M58").Value = -1413.7142
$ws.Range("N58").Value = -670788.0600000001
$ws.Range("H59").Value = 78701.664
$ws.Range("I59").Value = 55552.5
$ws.Range("J59").Value = 125000
$ws.Range("K59").Value = 55552.5
$ws.Range("L59").Value = 125000
$ws.Range("M59").Value = -54407.5
$ws.Range("N59").Value = -127290
$ws.Range("H99").Value = 2385.4517
$ws.Range("I99").Value = 2105.0952
$ws.Range("K99").Value = 2105.0952
$ws.Range("M99").Value = -607.0952000000002
$ws.Range("H107").Value = 1104.2354
$ws.Range("I107").Value = 1104.2354
$ws.Range("K107").Value = 1104.2354
$ws.Range("M107").Value = 815.7646
$ws.Range("H122").Value = 1996
$ws.Range("I122").Value = 1996
$ws.Range("K122").Value = 5988
$ws.Range("M122").Value = -3538
$ws.Range("H126").Value = 2385.4517
$ws.Range("I126").Value = 2105.0952
$ws.Range("K126").Value = 6315.285600000001
$ws.Range("M126").Value = -3845.285600000001
$ws.Range("H132").Value = 1112948
$ws.Range("I132").Value = 715576
$ws.Range("K132").Value = 2146728
$ws.Range("M132").Value = -2144198
$ws.Range("H134").Value = 2622.3547
$ws.Range("I134").Value = 2044.591
$ws.Range("K134").Value = 6133.772999999999
$ws.Range("M134").Value = -3598.772999999999
$ws.Range("H136").Value = 457593.1
$ws.Range("I136").Value = 1616.7142
$ws.Range("J136").Value = 670382.0600000001
$ws.Range("K136").Value = 4850.142599999999
$ws.Range("L136").Value = 2011146.18
$ws.Range("M136").Value = -2300.142599999999
$ws.Range("N136").Value = -2016246.18

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 605.5263
$ws.Range("I5").Value = 358.46155
$ws.Range("J5").Value = 1140.8334
$ws.Range("K5").Value = 1075.38465
$ws.Range("L5").Value = 3422.5002
$ws.Range("M5").Value = -963.38465
$ws.Range("N5").Value = -3646.5002
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H34").Value = 1782
$ws.Range("J34").Value = 2533.0715
$ws.Range("L34").Value = 7599.2145
$ws.Range("N34").Value = -7767.2145
$ws.Range("H68").Value = 1161.6666
$ws.Range("I68").Value = 1040
$ws.Range("K68").Value = 3120
$ws.Range("M68").Value = -2309
$ws.Range("H71").Value = 1161.6666
$ws.Range("I71").Value = 1040
$ws.Range("K71").Value = 9360
$ws.Range("M71").Value = -5304
$ws.Range("H92").Value = 2498.6667
$ws.Range("J92").Value = 2498.6667
$ws.Range("L92").Value = 7496.000100000001
$ws.Range("N92").Value = -9992.000100000001
$ws.Range("H119").Value = 10746.826
$ws.Range("I119").Value = 2764.125
$ws.Range("J119").Value = 15004.267
$ws.Range("K119").Value = 8292.375
$ws.Range("L119").Value = 45012.801
$ws.Range("M119").Value = -3454.375
$ws.Range("N119").Value = -54688.801
$ws.Range("H131").Value = 12477.35
$ws.Range("I131").Value = 850
$ws.Range("J131").Value = 13769.277
$ws.Range("K131").Value = 2550
$ws.Range("L131").Value = 41307.831
$ws.Range("M131").Value = 2490
$ws.Range("N131").Value = -51387.831
$ws.Range("H132").Value = 1634.625
$ws.Range("I132").Value = 1470
$ws.Range("J132").Value = 1909
$ws.Range("K132").Value = 13230
$ws.Range("L132").Value = 17181
$ws.Range("M132").Value = -10700
$ws.Range("N132").Value = -22241
$ws.Range("H135").Value = 605.5263
$ws.Range("I135").Value = 358.46155
$ws.Range("J135").Value = 1140.8334
$ws.Range("K135").Value = 3226.15395
$ws.Range("L135").Value = 10267.5006
$ws.Range("M135").Value = -691.1539499999999
$ws.Range("N135").Value = -15337.5006
$ws.Range("H138").Value = 3819998
$ws.Range("I138").Value = 5455542.5
$ws.Range("J138").Value = 221799.6
$ws.Range("K138").Value = 16366627.5
$ws.Range("L138").Value = 665398.8
$ws.Range("M138").Value = -16361487.5
$ws.Range("N138").Value = -675678.8
$ws.Range("H139").Value = 1854074.5
$ws.Range("I139").Value = 3032040.2
$ws.Range("K139").Value = 9096120.600000001
$ws.Range("M139").Value = -9090980.600000001
$ws.Range("H140").Value = 1110.6364
$ws.Range("I140").Value = 1110.6364
$ws.Range("K140").Value = 3331.9092
$ws.Range("M140").Value = 1848.0908

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H70").Value = 12890.583
$ws.Range("I70").Value = 13525.857
$ws.Range("J70").Value = 12001.2
$ws.Range("K70").Value = 13525.857
$ws.Range("L70").Value = 12001.2
$ws.Range("M70").Value = -13255.857
$ws.Range("N70").Value = -12541.2
$ws.Range("H73").Value = 12890.583
$ws.Range("I73").Value = 13525.857
$ws.Range("J73").Value = 12001.2
$ws.Range("K73").Value = 13525.857
$ws.Range("L73").Value = 12001.2
$ws.Range("M73").Value = -12589.857
$ws.Range("N73").Value = -13873.2
$ws.Range("H80").Value = 3305.158
$ws.Range("I80").Value = 3338.4443
$ws.Range("K80").Value = 3338.4443
$ws.Range("M80").Value = -2340.4443
$ws.Range("H83").Value = 3305.158
$ws.Range("I83").Value = 3338.4443
$ws.Range("K83").Value = 16692.2215
$ws.Range("M83").Value = -11700.2215
$ws.Range("H102").Value = 2353.6667
$ws.Range("I102").Value = 842.9
$ws.Range("J102").Value = 3727.0908
$ws.Range("K102").Value = 842.9
$ws.Range("L102").Value = 3727.0908
$ws.Range("M102").Value = 779.1
$ws.Range("N102").Value = -6971.0908
$ws.Range("H113").Value = 3000.28
$ws.Range("J113").Value = 3831.7856
$ws.Range("L113").Value = 3831.7856
$ws.Range("N113").Value = -8171.7856
$ws.Range("H122").Value = 55557908
$ws.Range("I122").Value = 2336.5
$ws.Range("K122").Value = 7009.5
$ws.Range("M122").Value = -4559.5
$ws.Range("H123").Value = 42354.332
$ws.Range("J123").Value = 42354.332
$ws.Range("L123").Value = 42354.332
$ws.Range("N123").Value = -47254.332
$ws.Range("H126").Value = 10304.056
$ws.Range("I126").Value = 19717.715
$ws.Range("K126").Value = 59153.145
$ws.Range("M126").Value = -56683.145
$ws.Range("H132").Value = 247773.92
$ws.Range("I132").Value = 360872.22
$ws.Range("K132").Value = 1082616.66
$ws.Range("M132").Value = -1080086.66
$ws.Range("H135").Value = 145159.8
$ws.Range("J135").Value = 154449.75
$ws.Range("L135").Value = 154449.75
$ws.Range("N135").Value = -164589.75
$ws.Range("H139").Value = 119222.164
$ws.Range("J139").Value = 119222.164
$ws.Range("L139").Value = 119222.164
$ws.Range("N139").Value = -129502.164

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 33336934
$ws.Range("I7").Value = 83336080
$ws.Range("J7").Value = 4166
$ws.Range("K7").Value = 83336080
$ws.Range("L7").Value = 4166
$ws.Range("M7").Value = -83335968
$ws.Range("N7").Value = -4390
$ws.Range("H14").Value = 7000
$ws.Range("J14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("N14").Value = -7344
$ws.Range("H22").Value = 3587174
$ws.Range("I22").Value = 3083.1667
$ws.Range("K22").Value = 3083.1667
$ws.Range("M22").Value = -2788.1667
$ws.Range("H27").Value = 3587174
$ws.Range("I27").Value = 3083.1667
$ws.Range("K27").Value = 3083.1667
$ws.Range("M27").Value = -2976.1667
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H32").Value = 3056.5
$ws.Range("I32").Value = 3056.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3056.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2739.5
$ws.Range("N32").ClearContents()
$ws.Range("H40").Value = 2822
$ws.Range("I40").Value = 2926.4546
$ws.Range("J40").Value = 2247.5
$ws.Range("K40").Value = 2926.4546
$ws.Range("L40").Value = 2247.5
$ws.Range("M40").Value = -2790.4546
$ws.Range("N40").Value = -2519.5
$ws.Range("H55").Value = 360.33334
$ws.Range("I55").Value = 395.42105
$ws.Range("J55").Value = 227
$ws.Range("K55").Value = 395.42105
$ws.Range("L55").Value = 227
$ws.Range("M55").Value = -222.42105
$ws.Range("N55").Value = -573
$ws.Range("H61").Value = 244228.92
$ws.Range("I61").Value = 283878.2
$ws.Range("J61").Value = 6333.3335
$ws.Range("K61").Value = 283878.2
$ws.Range("L61").Value = 6333.3335
$ws.Range("M61").Value = -283676.2
$ws.Range("N61").Value = -6737.3335
$ws.Range("H68").Value = 13312.125
$ws.Range("I68").Value = 2500
$ws.Range("J68").Value = 14856.714
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 14856.714
$ws.Range("M68").Value = -1751
$ws.Range("N68").Value = -16354.714
$ws.Range("H71").Value = 13312.125
$ws.Range("I71").Value = 2500
$ws.Range("J71").Value = 14856.714
$ws.Range("K71").Value = 12500
$ws.Range("L71").Value = 74283.57000000001
$ws.Range("M71").Value = -8756
$ws.Range("N71").Value = -81771.57000000001
$ws.Range("H93").Value = 5000
$ws.Range("I93").Value = 5000
$ws.Range("K93").Value = 5000
$ws.Range("M93").Value = -3752
$ws.Range("H113").Value = 244228.92
$ws.Range("I113").Value = 283878.2
$ws.Range("J113").Value = 6333.3335
$ws.Range("K113").Value = 283878.2
$ws.Range("L113").Value = 6333.3335
$ws.Range("M113").Value = -281708.2
$ws.Range("N113").Value = -10673.3335
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H126").Value = 33336934
$ws.Range("I126").Value = 83336080
$ws.Range("J126").Value = 4166
$ws.Range("K126").Value = 250008240
$ws.Range("L126").Value = 12498
$ws.Range("M126").Value = -250005770
$ws.Range("N126").Value = -17438
$ws.Range("H130").Value = 39849
$ws.Range("J130").Value = 39849
$ws.Range("L130").Value = 39849
$ws.Range("N130").Value = -49889
$ws.Range("H132").Value = 5270.7646
$ws.Range("I132").Value = 3131.5
$ws.Range("K132").Value = 9394.5
$ws.Range("M132").Value = -6864.5
$ws.Range("H136").Value = 4316.5273
$ws.Range("I136").Value = 3702.4888
$ws.Range("J136").Value = 7079.7
$ws.Range("K136").Value = 11107.4664
$ws.Range("L136").Value = 21239.1
$ws.Range("M136").Value = -8557.466400000001
$ws.Range("N136").Value = -26339.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 66913
$ws.Range("J46").Value = 66913
$ws.Range("L46").Value = 66913
$ws.Range("N46").Value = -67375
$ws.Range("H81").Value = 18189832
$ws.Range("I81").Value = 4400
$ws.Range("K81").Value = 8800
$ws.Range("M81").Value = -7739
$ws.Range("H84").Value = 18189832
$ws.Range("I84").Value = 4400
$ws.Range("K84").Value = 44000
$ws.Range("M84").Value = -38696
$ws.Range("H95").Value = 19990
$ws.Range("J95").Value = 19990
$ws.Range("L95").Value = 19990
$ws.Range("N95").Value = -25482
$ws.Range("H96").Value = 7011.4443
$ws.Range("I96").Value = 6652.25
$ws.Range("J96").Value = 7298.8
$ws.Range("K96").Value = 6652.25
$ws.Range("L96").Value = 7298.8
$ws.Range("M96").Value = -5279.25
$ws.Range("N96").Value = -10044.8
$ws.Range("H113").Value = 1131
$ws.Range("I113").Value = 1277.125
$ws.Range("K113").Value = 3831.375
$ws.Range("M113").Value = -1661.375
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 1848.0714
$ws.Range("I122").Value = 1709.3636
$ws.Range("K122").Value = 5128.0908
$ws.Range("M122").Value = -2678.0908
$ws.Range("H126").Value = 1436.5454
$ws.Range("I126").Value = 1475.25
$ws.Range("K126").Value = 4425.75
$ws.Range("M126").Value = -1955.75
$ws.Range("H132").Value = 492249.6
$ws.Range("I132").Value = 873785.2
$ws.Range("J132").Value = 4731.8887
$ws.Range("K132").Value = 2621355.6
$ws.Range("L132").Value = 14195.6661
$ws.Range("M132").Value = -2618825.6
$ws.Range("N132").Value = -19255.6661
$ws.Range("H134").Value = 66913
$ws.Range("J134").Value = 66913
$ws.Range("L134").Value = 200739
$ws.Range("N134").Value = -205809
$ws.Range("H136").Value = 3565.5173
$ws.Range("I136").Value = 2605.8333
$ws.Range("J136").Value = 5135.909
$ws.Range("K136").Value = 7817.499899999999
$ws.Range("L136").Value = 15407.727
$ws.Range("M136").Value = -5267.499899999999
$ws.Range("N136").Value = -20507.727
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 33736.332
$ws.Range("J140").Value = 33736.332
$ws.Range("L140").Value = 33736.332
$ws.Range("N140").Value = -44096.332
$ws.Range("H141").Value = 44999.4
$ws.Range("J141").Value = 44999.4
$ws.Range("L141").Value = 44999.4
$ws.Range("N141").Value = -55359.4
